# Apply the "Std Proton Afinity" column addition + header typo fix
# described by the commit "modified ml models with proton affinity".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the existing header typo: "Proton Affinity" -> "Proton Afinity"
$ws.Range("B1").Value = "Proton Afinity"

# 2) Add the new header in C1, matching the bold/centered/bordered style
#    already used by A1/B1.
$ws.Range("C1").Value = "Std Proton Afinity"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4160
$ws.Range("C1").Borders.LineStyle = 1

# 3) Populate the new "Std Proton Afinity" values for rows 2-43 (C2:C43).
$stdProtonAfinity = @(
    -1.007937091,
    -0.8047212419999999,
    1.643440602,
    -0.035446242,
    0.03311951,
    -0.562945989,
    -0.498655799,
    1.542646589,
    0.166551046,
    -1.328050909,
    -1.276673853,
    -0.664144091,
    0,
    0,
    0,
    2.238693251,
    2.35622384,
    -0.418593619,
    1.69111514,
    -0.117591686,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0.175367303,
    -0.194128568,
    0.222485708,
    -0.295803808,
    -0.312734515,
    -0.162355516,
    -0.461968928,
    -0.696767406,
    1.243828267,
    -0.683993846,
    -0.8549953729999999,
    -0.484176849,
    -0.451785926
)

for ($i = 0; $i -lt $stdProtonAfinity.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $stdProtonAfinity[$i]
}
